# Atualiza a aba "situacao_ordem_servicos" (lista de status da O.S.):
# reordena os status ABERTA / ATENDIMENTO EM ANDAMENTO / FINALIZADA / FATURADA / ...
# e torna essa aba a aba ativa (em vez de "servicos"), refletindo a mudança
# no fluxo de abertura/fechamento da O.S. descrita no commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("situacao_ordem_servicos")

$ws.Range("A3").Value = "ATENDIMENTO EM ANDAMENTO"
$ws.Range("A4").Value = "FINALIZADA"
$ws.Range("A5").Value = "FATURADA"

$ws.Activate()
$ws.Range("A5").Select()
